# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off (new handoff xliff files generated, status updated, and an
# error detail noting the handback file is stale).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8db7f4ba2bf9ba57cdd2327cc6107eefe289827f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa7a424c168658858a944536ebdeae2c7648dd5f/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 18:48:33"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("P1").ColumnWidth = 39.17
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-20 18:48:29"
$wsZhCn.Range("P3").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("P1").ColumnWidth = 39.17
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-20 18:48:33"
$wsDeDe.Range("P3").Value = $errorDetail

Write-Host "Handoff report generated."
